$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.786.12"
$ws.Range("E2").Value = "  +0.70%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.701.65"
$ws.Range("E3").Value = "  +0.52%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.59"
$ws.Range("E5").Value = "  +0.05%  "

# Row 6
$ws.Range("E6").Value = "  -0.11%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3983"
$ws.Range("E7").Value = "  +2.22%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4053"
$ws.Range("E8").Value = "  +0.35%  "

# Row 9
$ws.Range("B9").Value = "Polygon"
$ws.Range("C9").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.475"
$ws.Range("E9").Value = "  -1.50%  "

# Row 10
$ws.Range("B10").Value = "BinanceUSD"
$ws.Range("C10").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9998"
$ws.Range("E10").Value = "  -0.31%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.48"
$ws.Range("E11").Value = "  +1.05%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08791"
$ws.Range("E12").Value = "  +0.42%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.24"
$ws.Range("E13").Value = "  +3.01%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.500"
$ws.Range("E14").Value = "  -0.47%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.948"
$ws.Range("E15").Value = "  -0.18%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001340"
$ws.Range("E16").Value = "  -1.10%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.816.35"
$ws.Range("E17").Value = "  +7.66%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.45"
$ws.Range("E18").Value = "  -3.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07174"
$ws.Range("E19").Value = "  +0.85%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.89"
$ws.Range("E20").Value = "  +4.64%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.312"
$ws.Range("E21").Value = "  +0.21%  "

# Row 22
$ws.Range("E22").Value = "  -0.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.39"
$ws.Range("E23").Value = "  +0.81%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.780.81"
$ws.Range("E24").Value = "  +0.72%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.372"
$ws.Range("E25").Value = "  +0.80%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.881"
$ws.Range("E26").Value = "  -4.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.09"
$ws.Range("E27").Value = "  +1.16%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.052"
$ws.Range("E28").Value = "  +16.02%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.87"
$ws.Range("E29").Value = "  -0.24%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "143.78"
$ws.Range("E30").Value = "  +4.79%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.219"
$ws.Range("E31").Value = "  -6.24%  "

# Row 32
$ws.Range("B32").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.981.32"
$ws.Range("E32").Value = "  +5.77%  "

# Row 33
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.263"
$ws.Range("E33").Value = "  +14.71%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.03183"
$ws.Range("E34").Value = "  +8.31%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08573"
$ws.Range("E35").Value = "  -3.08%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.251"
$ws.Range("E36").Value = "  -2.15%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.028"
$ws.Range("E37").Value = "  -0.91%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2847"
$ws.Range("E38").Value = "  +3.45%  "

# Row 39
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09463"
$ws.Range("E39").Value = "  +3.45%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8325"
$ws.Range("E40").Value = "  +5.30%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.72"
$ws.Range("E41").Value = "  -0.85%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.18"
$ws.Range("E42").Value = "  -0.86%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.476"
$ws.Range("E43").Value = "  +0.94%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.46"
$ws.Range("E44").Value = "  +4.22%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.702"
$ws.Range("E45").Value = "  +4.77%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7402"
$ws.Range("E46").Value = "  +2.45%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.218"
$ws.Range("E47").Value = "  +0.20%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.375"
$ws.Range("E48").Value = "  +2.73%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.002"
$ws.Range("E49").Value = "  +0.04%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.08389"
$ws.Range("E50").Value = "  +5.25%  "

# Row 51
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "139.39"
$ws.Range("E51").Value = "  +1.06%  "
